$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Puliendo la estructura de los data frames de CPU y MEM
# - Columna B (CANTIDAD) homogeneizada a 29 en todas las filas de datos (2-13)
# - Columna C (MAX) actualizada en varias filas
# - Columna D (MIN) igualada al valor final de la columna C (MAX) de cada fila

function Set-TextValue($cell, [string]$text) {
    # Force the value to be stored as literal text (avoid Excel auto-converting
    # percentage-looking strings like "0.97%" into numeric percent values),
    # then reset the cell style so no extra number-format style lingers.
    $cell.Formula = "'" + $text
    $cell.Style = "Normal"
}

# Column B -> 29 for every data row
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 2).Value2 = 29
}

# Column C (MAX) updates
Set-TextValue $ws.Cells.Item(3, 3) "0.97%"
Set-TextValue $ws.Cells.Item(4, 3) "0.31%"
Set-TextValue $ws.Cells.Item(5, 3) "0.47%"
Set-TextValue $ws.Cells.Item(6, 3) "16.27%"
Set-TextValue $ws.Cells.Item(8, 3) "1.96%"
Set-TextValue $ws.Cells.Item(9, 3) "1.94%"
Set-TextValue $ws.Cells.Item(10, 3) "1.59%"

# Column D (MIN) -> set equal to column C (MAX) value for each row
for ($row = 2; $row -le 13; $row++) {
    $maxText = $ws.Cells.Item($row, 3).Text
    Set-TextValue $ws.Cells.Item($row, 4) $maxText
}
